$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 187, shifting existing rows 187-262 down to 188-263.
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(187, 1).Value = 11
$ws.Cells.Item(187, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(187, 3).Value = "Bíobío"
$ws.Cells.Item(187, 4).Value = 45006
$ws.Cells.Item(187, 5).Value = 8
$ws.Cells.Item(187, 6).Value = 100112003
$ws.Cells.Item(187, 7).Value = "Ajo"
$ws.Cells.Item(187, 8).Value = "Chino"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 300
$ws.Cells.Item(187, 11).Value = 15000
$ws.Cells.Item(187, 12).Value = 16000
$ws.Cells.Item(187, 13).Value = 15500
$ws.Cells.Item(187, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(187, 15).Value = "China"
$ws.Cells.Item(187, 16).Value = 1550
$ws.Cells.Item(187, 17).Value = 10
$ws.Cells.Item(187, 18).Value = "Hortaliza"
